$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 97.57717366666668
$ws.Range("H2").Value = 292.731521
$ws.Range("I2").Value = 0.3532166605548384
$ws.Range("J2").Value = 0.3532166605548384
$ws.Range("M2").Value = 13.11484166666667
$ws.Range("N2").Value = 39.344525
$ws.Range("O2").Value = 0.2389043281085165
$ws.Range("P2").Value = 0.2389043281085165
$ws.Range("Q2").Value = 1279.70918291917
$ws.Range("R2").Value = 11517.38264627253
$ws.Range("S2").Value = 0.08438498896658761
$ws.Range("T2").Value = 0.08438498896658762
$ws.Range("G3").Value = 97.57717366666668
$ws.Range("H3").Value = 292.731521
$ws.Range("I3").Value = 0.3532166605548384
$ws.Range("J3").Value = 0.3532166605548384
$ws.Range("O3").Value = 0.5367702700792449
$ws.Range("P3").Value = 0.5367702700792449
$ws.Range("Q3").Value = 2875.250729766605
$ws.Range("R3").Value = 25877.25656789945
$ws.Range("S3").Value = 0.1895962022825096
$ws.Range("T3").Value = 0.1895962022825096
$ws.Range("G4").Value = 97.57717366666668
$ws.Range("H4").Value = 292.731521
$ws.Range("I4").Value = 0.3532166605548384
$ws.Range("J4").Value = 0.3532166605548384
$ws.Range("O4").Value = 0.2243254018122386
$ws.Range("P4").Value = 0.2243254018122386
$ws.Range("Q4").Value = 1201.61605666164
$ws.Range("R4").Value = 10814.54450995476
$ws.Range("S4").Value = 0.0792354693057412
$ws.Range("T4").Value = 0.07923546930574123
$ws.Range("I5").Value = 0.5533024543641269
$ws.Range("J5").Value = 0.5533024543641269
$ws.Range("M5").Value = 13.11484166666667
$ws.Range("N5").Value = 39.344525
$ws.Range("O5").Value = 0.2389043281085165
$ws.Range("P5").Value = 0.2389043281085165
$ws.Range("Q5").Value = 2004.622971830509
$ws.Range("R5").Value = 18041.60674647458
$ws.Range("S5").Value = 0.1321863511006549
$ws.Range("T5").Value = 0.1321863511006549
$ws.Range("I6").Value = 0.5533024543641269
$ws.Range("J6").Value = 0.5533024543641269
$ws.Range("O6").Value = 0.5367702700792449
$ws.Range("P6").Value = 0.5367702700792449
$ws.Range("S6").Value = 0.2969963078645415
$ws.Range("T6").Value = 0.2969963078645415
$ws.Range("I7").Value = 0.5533024543641269
$ws.Range("J7").Value = 0.5533024543641269
$ws.Range("O7").Value = 0.2243254018122386
$ws.Range("P7").Value = 0.2243254018122386
$ws.Range("S7").Value = 0.1241197953989306
$ws.Range("T7").Value = 0.1241197953989306
$ws.Range("I8").Value = 0.09348088508103472
$ws.Range("J8").Value = 0.09348088508103473
$ws.Range("M8").Value = 13.11484166666667
$ws.Range("N8").Value = 39.344525
$ws.Range("O8").Value = 0.2389043281085165
$ws.Range("P8").Value = 0.2389043281085165
$ws.Range("Q8").Value = 338.6826286101501
$ws.Range("R8").Value = 3048.14365749135
$ws.Range("S8").Value = 0.02233298804127405
$ws.Range("T8").Value = 0.02233298804127405
$ws.Range("I9").Value = 0.09348088508103472
$ws.Range("J9").Value = 0.09348088508103473
$ws.Range("O9").Value = 0.5367702700792449
$ws.Range("P9").Value = 0.5367702700792449
$ws.Range("S9").Value = 0.05017775993219387
$ws.Range("T9").Value = 0.05017775993219387
$ws.Range("I10").Value = 0.09348088508103472
$ws.Range("J10").Value = 0.09348088508103473
$ws.Range("O10").Value = 0.2243254018122386
$ws.Range("P10").Value = 0.2243254018122386
$ws.Range("S10").Value = 0.02097013710756681
$ws.Range("T10").Value = 0.02097013710756682
